# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 0.08108108108108109
$ws.Range("D3").Value = 0.3153153153153153
$ws.Range("E3").Value = 0.7027027027027027
$ws.Range("F3").Value = 0.9279279279279279
$ws.Range("H3").Value = 0.1296509397775221
$ws.Range("I3").Value = 0.385381854517657
$ws.Range("J3").Value = -0.01801801801801802
$ws.Range("K3").Value = 2022.603603603604

$ws.Range("Q3").Value = 524
$ws.Range("R3").Value = 987
$ws.Range("S3").Value = 1836
$ws.Range("T3").Value = 2817
$ws.Range("U3").Value = 3583
$ws.Range("V3").Value = 4579
$ws.Range("W3").Value = 4116
$ws.Range("X3").Value = 3267
$ws.Range("Y3").Value = 2286
$ws.Range("Z3").Value = 1520

$ws.Range("AF3").Value = 0.897315
$ws.Range("AG3").Value = 0.806584
$ws.Range("AH3").Value = 0.640212
$ws.Range("AI3").Value = 0.447972
$ws.Range("AJ3").Value = 0.297864
